$p = $ppt.ActivePresentation

# The agenda slide is slide 5 ("1 Introduction.pptx" deck), shape "TextShape 2"
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item("TextShape 2")
$tf = $sh.TextFrame
$tr = $tf.TextRange

# ---------------------------------------------------------------------------
# 1. Remove the old "Exploitation Related Activities (WP7)" bullet that used
#    to sit right after the Break bullet (paragraph 7) - the agenda moved it
#    earlier in the schedule.
# ---------------------------------------------------------------------------
$tr.Paragraphs(7, 1).Delete()

# ---------------------------------------------------------------------------
# 2. Break bullet: 14:35-14:50 -> 14:40-15:00
# ---------------------------------------------------------------------------
$tr.Paragraphs(6, 1).Runs(1, 1).Text = "Break (15 mins) [14:40 - 15:00]"

# ---------------------------------------------------------------------------
# 3. User Interfaces & Integration (WP6): 14:25-14:35 -> 14:30-14:40
# ---------------------------------------------------------------------------
$tr.Paragraphs(5, 1).Runs(2, 1).Text = "(10 mins) [14:30 - 14:40]"

# ---------------------------------------------------------------------------
# 4. Technical Presentations (WP2, 3 & 4): "(45 mins) [13:40 - 14:25]"
#    becomes three runs: "(" + "30" + " mins) [14:00 - 14:30]"
# ---------------------------------------------------------------------------
$para4 = $tr.Paragraphs(4, 1)
$timeRun = $para4.Runs(2, 1)
$timeRun.Text = "("
$null = $timeRun.InsertAfter("30")
$null = $para4.Runs(3, 1).InsertAfter(" mins) [14:00 - 14:30]")

# ---------------------------------------------------------------------------
# 5. Industrial Use Cases (WP1): 20 min [13:20 - 13:40] -> 30 min [13:20 - 13:50]
# ---------------------------------------------------------------------------
$para3 = $tr.Paragraphs(3, 1)
$para3.Runs(2, 1).Text = "(30 min) [13:20 - 13:50]"

# ---------------------------------------------------------------------------
# 6. Insert a brand-new bullet right after "Industrial Use Cases" for the
#    "Exploitation Related Activities (WP7)" item that was moved up in the
#    agenda, now scheduled [13:50 - 14:00].
# ---------------------------------------------------------------------------
$lastRun = $para3.Runs($para3.Runs().Count, 1)
$lastRun.InsertAfter("`rExploitation Related Activities (WP7) (10 mins) [13:50 - 14:00]") | Out-Null

$newPara = $tr.Paragraphs(4, 1)
$firstLen = "Exploitation Related Activities (WP7) ".Length
$titlePart = $newPara.Characters(1, $firstLen)
$titlePart.Font.Bold = $false
$titlePart.Font.Color.RGB = 255
